$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row just above the "state (Rust)" row (row 29) so the new
# "IconPark" component lands at row 29 and everything below shifts down
# by one (state (Rust) -> 30, Microsoft Visual Studio Runtimes -> 31).
$ws.Rows.Item(29).Insert()

# Copy the formatting of the row above (row 28) onto the freshly inserted
# row so the new row reuses the existing "data row" cell styles instead of
# generating brand-new style entries.
$ws.Range("B28:H28").Copy()
$ws.Range("B29:H29").PasteSpecial(-4122)

# Fill in the new component's data.
$ws.Range("C29").Value() = "IconPark"
$ws.Range("D29").Value() = "any"
$ws.Range("E29").Value() = "Apache v 2.0"
$ws.Range("F29").Value() = "github.com/bytedance/IconPark"
$ws.Range("G29").Value() = "runtime"
$ws.Range("H29").Value() = "no"

# The merged "Third Party Licenses" cell in column B must grow to keep
# covering the whole (now one-row-taller) table.
$ws.Range("B4:B31").Merge()

# Restore the selection that was active when the file was last saved.
$ws.Range("F9").Select()
